$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.173.61'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.900.02'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = '''252.63'
$ws.Range('E5').Value = '  +2.80%  '
$ws.Range('D6').Value = '''0.698'
$ws.Range('E6').Value = '  +1.22%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').Value = '''40.89'
$ws.Range('E8').Value = '  -2.37%  '
$ws.Range('D9').Value = '''0.361'
$ws.Range('E9').Value = '  +3.34%  '
$ws.Range('D10').Value = '''52.87'
$ws.Range('E10').Value = '  -0.30%  '
$ws.Range('D11').Value = '''0.0752'
$ws.Range('E11').Value = '  +3.69%  '
$ws.Range('D12').Value = '''0.0982'
$ws.Range('E12').Value = '  -1.23%  '
$ws.Range('D13').Value = '''13.08'
$ws.Range('E13').Value = '  +6.24%  '
$ws.Range('D14').Value = '2.174.73'
$ws.Range('E15').Value = '  +3.43%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.940.56'
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '''4.95'
$ws.Range('E17').Value = '  +2.28%  '
$ws.Range('D18').Value = '35.174.52'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').Value = '''73.66'
$ws.Range('E19').Value = '  +1.69%  '
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('D21').Value = '''242.18'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('D22').Value = '''12.94'
$ws.Range('E22').Value = '  +2.74%  '
$ws.Range('D23').Value = '''5.04'
$ws.Range('E23').Value = '  +4.22%  '
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('E25').Value = '  +3.85%  '
$ws.Range('D26').Value = '''2.27'
$ws.Range('E26').Value = '  -1.53%  '
$ws.Range('D27').Value = '''166.80'
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('D28').Value = '''8.58'
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('D29').Value = '''18.47'
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('D30').Value = '''0.130'
$ws.Range('E30').Value = '  -0.71%  '
$ws.Range('D31').Value = '4.128.88'
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').Value = '''2.07'
$ws.Range('E32').Value = '  +15.86%  '
$ws.Range('D33').Value = '''0.0603'
$ws.Range('E33').Value = '  +5.55%  '
$ws.Range('D34').Value = '''4.32'
$ws.Range('E34').Value = '  +3.43%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '''4.20'
$ws.Range('E35').Value = '  +1.77%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').Value = '''1.55'
$ws.Range('E36').Value = '  +15.94%  '
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('D38').Value = '''0.855'
$ws.Range('E38').Value = '  -11.27%  '
$ws.Range('D39').Value = '''2.01'
$ws.Range('E39').Value = '  -1.33%  '
$ws.Range('D40').Value = '''100.66'
$ws.Range('E40').Value = '  +11.26%  '
$ws.Range('D41').Value = '''17.14'
$ws.Range('E41').Value = '  +6.03%  '
$ws.Range('D42').Value = '''0.0214'
$ws.Range('E42').Value = '  +1.77%  '
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').Value = '''0.0648'
$ws.Range('E44').Value = '  -5.35%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '''2.43'
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.313.94'
$ws.Range('E46').Value = '  -2.34%  '
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('E48').Value = '  -1.65%  '
$ws.Range('D49').Value = '''6.58'
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('D50').Value = '''11.79'
$ws.Range('E50').Value = '  -6.41%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = '''43.10'
$ws.Range('E51').Value = '  -8.46%  '
